# Update Name of Algo
# Applies the revised RandomForest imputation values for the
# terrestrial_mammals / combination_2_ABCDE / BD / 20 / seed5 run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    "D3"  = -7.137099999999995
    "D4"  = -6.892899999999994
    "B11" = 5.9689
    "B12" = 4.969399999999998
    "D14" = -7.4749
    "B15" = 5.313299999999998
    "D26" = -8.45820000000001
    "B27" = 6.175700000000005
    "B28" = 6.112900000000003
    "B31" = 6.165800000000001
    "D31" = -8.325899999999999
    "B32" = 6.5665
    "D35" = -8.386099999999994
    "B36" = 8.9877
    "D37" = -7.476499999999998
    "B38" = 5.2648
    "D39" = -8.056500000000003
    "D40" = -8.014699999999996
    "D45" = -7.849599999999998
    "B46" = 6.8546
    "D52" = -7.605499999999997
    "B54" = 4.830399999999998
    "B55" = 5.5198
    "B56" = 5.1605
    "D57" = -8.4214
    "B67" = 5.438699999999995
    "B69" = 5.233499999999997
    "B72" = 5.158600000000002
    "B73" = 8.387700000000001
    "D81" = -7.065399999999991
    "B83" = 5.237399999999997
    "D83" = -8.909300000000002
    "B86" = 4.944800000000002
    "B91" = 5.364500000000001
    "B93" = 5.709800000000001
    "B99" = 4.651899999999995
    "D100" = -8.057700000000001
    "D102" = -7.478099999999996
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
